$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression
$ws.Range("B2").Value = 2920447597530898

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 8463869133625.325
$ws.Range("C3").Value = 7571230317425.839
$ws.Range("D3").Value = 6848878418676.227

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2975774401670.688
$ws.Range("C4").Value = 2819154696319.6
$ws.Range("D4").Value = 3054084254346.23

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 90146323950372.44
$ws.Range("C5").Value = 141214467684338.8
$ws.Range("D5").Value = 197909518742438
